$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @{ A = 90705; B = "Clarice Mendes";          C = "P&D";               D = "Problemas pessoais"; E = 1; F = 45100; G = 4875.92 }
  3  = @{ A = 76324; B = "João Vitor Almeida";      C = "Recursos Humanos";  D = "Consulta medica";    E = 8; F = 45102; G = 2047.7  }
  4  = @{ A = 56110; B = "Sr. Liam da Rosa";        C = "Financeiro";        D = "Doenca";              E = 4; F = 45083; G = 5223.31 }
  5  = @{ A = 92763; B = "Sra. Juliana Duarte";     C = "Engenharia";        D = "Outros";              E = 8; F = 45092; G = 4654.48 }
  6  = @{ A = 3721;  B = "Henry Gabriel Oliveira";  C = "Engenharia";        D = "Doenca";              E = 6; F = 45105; G = 8479.5  }
  7  = @{ A = 96593; B = "Vitor Cassiano";          C = "Marketing";         D = "Consulta medica";    E = 2; F = 45085; G = 5221.23 }
  8  = @{ A = 31726; B = "Agatha Costela";          C = "P&D";               D = "Outros";              E = 7; F = 45082; G = 2932.38 }
  9  = @{ A = 75883; B = "João Pedro Nunes";        C = "Marketing";         D = "Problemas pessoais"; E = 8; F = 45081; G = 6024.03 }
  10 = @{ A = 11838; B = "Mariah Costa";            C = "Marketing";         D = "Doenca";              E = 1; F = 45105; G = 3564.38 }
  11 = @{ A = 29132; B = "Calebe Duarte";           C = "Vendas";            D = "Outros";              E = 5; F = 45095; G = 4633.74 }
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Cells.Item($row, 1).Value = $vals.A
  $ws.Cells.Item($row, 2).Value = $vals.B
  $ws.Cells.Item($row, 3).Value = $vals.C
  $ws.Cells.Item($row, 4).Value = $vals.D
  $ws.Cells.Item($row, 5).Value = $vals.E
  $ws.Cells.Item($row, 6).Value = $vals.F
  $ws.Cells.Item($row, 7).Value = $vals.G
}
